$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kenntnisse")
$ws.Activate()

# Update the quantitative value for Französisch (row 11, column C) from 2 to 1
$ws.Range("C11").Value = 1

# Update the active selection on the sheet to C15
$ws.Range("C15").Select()
